$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.734.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.701.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.33%  "
$ws.Range("E6").Value = "  +0.28%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3942"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.34%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "54.13"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("B11").Value = "BinanceUSD"
$ws.Range("C11").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.005"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08877"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.11%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.472"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.68"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.202"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +7.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001325"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.709.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.86%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07061"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.077"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.005"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.75"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.69%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.718.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.51%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.152"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.25%  "
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.77"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "162.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.702"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +16.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.92"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.174"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.23%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09087"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.677"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.074"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.981"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.63%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2760"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02783"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09159"
$ws.Range("D40").Style = "Normal"
$ws.Range("E41").Value = "  -0.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.7703"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.45%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7189"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.580"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.220"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("E47").Value = "  +0.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.346"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.30%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "91.14"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07990"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.66%  "
